# Applies the "ifoCAST full series evaluation" update to the forecast
# error table: the evaluation window advanced by one quarter, so the
# ME/MAE/MSE/RMSE/SE/N statistics for quarters Q0..Q9 (rows 2..11) are
# refreshed with recomputed values, and row 11 (Q9) now also has an SE
# value in column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.1420658796491534, 0.2963740997880928, 0.1213439170623603, 0.3483445378678418, 0.3310457719767695, 13)
    3  = @(0.1568915331285771, 0.2900446020859644, 0.1128329969643262, 0.3359062323987547, 0.3102222847576306, 12)
    4  = @(0.1956053612568337, 0.2558683862338031, 0.08648424928328589, 0.2940820451562555, 0.2303151560882531, 11)
    5  = @(0.222511934274951, 0.2817985532879059, 0.1129743571721927, 0.3361165827093223, 0.2655451338021729, 10)
    6  = @(0.176450614684246, 0.1849112350951672, 0.06089664821255769, 0.2467724624275523, 0.1829810301339271, 9)
    7  = @(0.2112811930113801, 0.285830133390834, 0.1115765264066959, 0.3340307267403643, 0.2765848542209762, 8)
    8  = @(0.1601567932497119, 0.2716627869588349, 0.08940684574059213, 0.2990097753261457, 0.2727320208605046, 7)
    9  = @(0.2627214152773011, 0.3235434349789264, 0.1669924024756296, 0.4086470389904099, 0.3380801243818421, 7)
    10 = @(0.07165448440387559, 0.1716308632352222, 0.0355562541576917, 0.1885636607559678, 0.1950060544653244, 5)
    11 = @(0.08333287966522271, 0.2888136092683156, 0.1119970241814113, 0.3346595646047059, 0.3742595451432296, 4)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]  # B - ME
    $ws.Cells.Item($row, 3).Value = $values[1]  # C - MAE
    $ws.Cells.Item($row, 4).Value = $values[2]  # D - MSE
    $ws.Cells.Item($row, 5).Value = $values[3]  # E - RMSE
    $ws.Cells.Item($row, 6).Value = $values[4]  # F - SE
    $ws.Cells.Item($row, 7).Value = $values[5]  # G - N
}
